# Bitácora: add the "RoundTable" entry (Mesa redonda) to the
# "Actividades realizadas" log, as described in the commit message:
#   "Se agregó en la bitácora la información sobre los tiempos de la
#    mesa redonda."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actividades realizadas")

# The table currently ends at row 14. Copy the formatting of the last
# data row down into the new row 15 so the new entry keeps the same
# look (fonts/number format/borders) as the rest of the log.
$ws.Range("B14:I14").Copy()
$ws.Range("B15:I15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New bitácora entry: RoundTable activity.
$ws.Range("B15").Value = 43512
$ws.Range("C15").Value = "Mauricio"
$ws.Range("D15").Value = "Medium"
$ws.Range("E15").Value = "App requests"
$ws.Range("F15").Value = "RoundTable"
$ws.Range("G15").Value = "1hr 30 min (19:00 - 20:30)"
$ws.Range("H15").Value = "Great, it was interesting to gather new points of view and the activity was enjoyable."
$ws.Range("I15").Value = "Some other reasons that causes people to feel depressed."

# Leave the selection where the author left it after typing the new row.
$ws.Range("I17").Select()
